# get ngo profile details Api added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 13: NGO profile details endpoint (set in this order so the
# shared-string table grows in the same sequence as the target workbook)
$ws.Range("A13").Value = "/api/ngo/profiledetails/:userid"
$ws.Range("B13").Value = "get"
$ws.Range("C13").Value = "get all details of ngo"

# Reword the existing "user profile details" row's description
$ws.Range("C12").Value = "get all details user"

# Match the formatting used by the other data rows (bigger font on column A,
# 15.75pt row height)
$ws.Range("A13").Font.Size = 12
$ws.Rows.Item(13).RowHeight = 15.75

$ws.Range("A14").Select()
